# Auto-generated edit script applying the Balmung_Profits.xlsx value updates
# across the ALC, ARM, BSM, CRP, CUL, GSM and LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 90912424
$ws.Range("I106").Value = 142860080
$ws.Range("J106").Value = 4024.75
$ws.Range("K106").Value = 142860080
$ws.Range("L106").Value = 4024.75
$ws.Range("M106").Value = -142859449
$ws.Range("N106").Value = -5286.75
$ws.Range("H137").Value = 6669532
$ws.Range("I137").Value = 1862.5
$ws.Range("K137").Value = 5587.5
$ws.Range("M137").Value = -3037.5
$ws.Range("H138").Value = 6117.7144
$ws.Range("I138").Value = 10660.723
$ws.Range("J138").Value = 4300.511
$ws.Range("K138").Value = 31982.169
$ws.Range("L138").Value = 12901.533
$ws.Range("M138").Value = -26842.169
$ws.Range("N138").Value = -23181.533
$ws.Range("H141").Value = 1768.359
$ws.Range("I141").Value = 1294.0646
$ws.Range("K141").Value = 3882.1938
$ws.Range("M141").Value = 1297.8062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 739.625
$ws.Range("I4").Value = 723.8570999999999
$ws.Range("K4").Value = 723.8570999999999
$ws.Range("M4").Value = -607.8570999999999
$ws.Range("H6").Value = 10000
$ws.Range("I6").Value = 10000
$ws.Range("K6").Value = 10000
$ws.Range("M6").Value = -9827
$ws.Range("H32").Value = 106463.914
$ws.Range("I32").Value = 109576.78
$ws.Range("K32").Value = 109576.78
$ws.Range("M32").Value = -109289.78
$ws.Range("H61").Value = 1054697.5
$ws.Range("I61").Value = 3716.7727
$ws.Range("K61").Value = 3716.7727
$ws.Range("M61").Value = -3504.7727
$ws.Range("H74").Value = 1175217.8
$ws.Range("I74").Value = 2378
$ws.Range("J74").Value = 1487975.1
$ws.Range("K74").Value = 2378
$ws.Range("L74").Value = 1487975.1
$ws.Range("M74").Value = -1504
$ws.Range("N74").Value = -1489723.1
$ws.Range("H77").Value = 1175217.8
$ws.Range("I77").Value = 2378
$ws.Range("J77").Value = 1487975.1
$ws.Range("K77").Value = 11890
$ws.Range("L77").Value = 7439875.5
$ws.Range("M77").Value = -7522
$ws.Range("N77").Value = -7448611.5
$ws.Range("H119").Value = 63538.57
$ws.Range("J119").Value = 63538.57
$ws.Range("L119").Value = 63538.57
$ws.Range("N119").Value = -73214.57000000001
$ws.Range("H124").Value = 38107.5
$ws.Range("J124").Value = 38107.5
$ws.Range("L124").Value = 38107.5
$ws.Range("N124").Value = -47927.5
$ws.Range("H125").Value = 82500
$ws.Range("J125").Value = 82500
$ws.Range("L125").Value = 82500
$ws.Range("N125").Value = -92340
$ws.Range("H128").Value = 95370.336
$ws.Range("J128").Value = 95370.336
$ws.Range("L128").Value = 95370.336
$ws.Range("N128").Value = -105330.336
$ws.Range("H130").Value = 49666.668
$ws.Range("J130").Value = 49666.668
$ws.Range("L130").Value = 49666.668
$ws.Range("N130").Value = -59706.668
$ws.Range("H131").Value = 60000
$ws.Range("J131").Value = 60000
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080
$ws.Range("H132").Value = 2344.5
$ws.Range("I132").Value = 1381.762
$ws.Range("J132").Value = 3899.6924
$ws.Range("K132").Value = 4145.286
$ws.Range("L132").Value = 11699.0772
$ws.Range("M132").Value = -1615.286
$ws.Range("N132").Value = -16759.0772
$ws.Range("H136").Value = 1054697.5
$ws.Range("I136").Value = 3716.7727
$ws.Range("K136").Value = 11150.3181
$ws.Range("M136").Value = -8600.3181
$ws.Range("H140").Value = 86068.86
$ws.Range("J140").Value = 86068.86
$ws.Range("L140").Value = 86068.86
$ws.Range("N140").Value = -96428.86
$ws.Range("H141").Value = 108994.43
$ws.Range("J141").Value = 108994.43
$ws.Range("L141").Value = 108994.43
$ws.Range("N141").Value = -119354.43

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 24000
$ws.Range("J62").Value = 24000
$ws.Range("L62").Value = 24000
$ws.Range("N62").Value = -25372
$ws.Range("H64").Value = 2565.4285
$ws.Range("I64").Value = 3012
$ws.Range("J64").Value = 2386.8
$ws.Range("K64").Value = 3012
$ws.Range("L64").Value = 2386.8
$ws.Range("M64").Value = -2787
$ws.Range("N64").Value = -2836.8
$ws.Range("H65").Value = 24000
$ws.Range("J65").Value = 24000
$ws.Range("L65").Value = 72000
$ws.Range("N65").Value = -78864
$ws.Range("H67").Value = 2565.4285
$ws.Range("I67").Value = 3012
$ws.Range("J67").Value = 2386.8
$ws.Range("K67").Value = 3012
$ws.Range("L67").Value = 2386.8
$ws.Range("M67").Value = -2232
$ws.Range("N67").Value = -3946.8
$ws.Range("H126").Value = 24593.334
$ws.Range("J126").Value = 24593.334
$ws.Range("L126").Value = 24593.334
$ws.Range("N126").Value = -34473.334
$ws.Range("H127").Value = 20000
$ws.Range("J127").Value = 20000
$ws.Range("L127").Value = 20000
$ws.Range("N127").Value = -29920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 42.833332
$ws.Range("I7").Value = 40.5
$ws.Range("J7").Value = 44
$ws.Range("K7").Value = 40.5
$ws.Range("L7").Value = 44
$ws.Range("M7").Value = 72.5
$ws.Range("N7").Value = -270
$ws.Range("H31").Value = 5192.2964
$ws.Range("I31").Value = 3336.5
$ws.Range("J31").Value = 5973.684
$ws.Range("K31").Value = 3336.5
$ws.Range("L31").Value = 5973.684
$ws.Range("M31").Value = -3041.5
$ws.Range("N31").Value = -6563.684
$ws.Range("H34").Value = 5192.2964
$ws.Range("I34").Value = 3336.5
$ws.Range("J34").Value = 5973.684
$ws.Range("K34").Value = 3336.5
$ws.Range("L34").Value = 5973.684
$ws.Range("M34").Value = -3134.5
$ws.Range("N34").Value = -6377.684
$ws.Range("H132").Value = 3065.1353
$ws.Range("I132").Value = 2944.5334
$ws.Range("J132").Value = 3582
$ws.Range("K132").Value = 8833.600199999999
$ws.Range("L132").Value = 10746
$ws.Range("M132").Value = -6303.600199999999
$ws.Range("N132").Value = -15806

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3354.1428
$ws.Range("I123").Value = 3354.1428
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 10062.4284
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -7612.428400000001
$ws.Range("N123").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 715690.4399999999
$ws.Range("I132").Value = 3090.76
$ws.Range("J132").Value = 2335235.2
$ws.Range("K132").Value = 9272.280000000001
$ws.Range("L132").Value = 7005705.600000001
$ws.Range("M132").Value = -6742.280000000001
$ws.Range("N132").Value = -7010765.600000001
$ws.Range("H137").Value = 15780
$ws.Range("J137").Value = 15780
$ws.Range("L137").Value = 15780
$ws.Range("N137").Value = -25980
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = $null
$ws.Range("N138").Value = $null
$ws.Range("H139").Value = 374999.66
$ws.Range("J139").Value = 374999.66
$ws.Range("L139").Value = 374999.66
$ws.Range("N139").Value = -385279.66
$ws.Range("H140").Value = 199854.5
$ws.Range("I140").Value = 199709
$ws.Range("J140").Value = 200000
$ws.Range("K140").Value = 199709
$ws.Range("L140").Value = 200000
$ws.Range("M140").Value = -194529
$ws.Range("N140").Value = -210360
$ws.Range("H141").Value = 62764.223
$ws.Range("J141").Value = 62764.223
$ws.Range("L141").Value = 62764.223
$ws.Range("N141").Value = -73124.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15580
$ws.Range("H68").Value = 3831.4666
$ws.Range("I68").Value = 3719.2222
$ws.Range("K68").Value = 3719.2222
$ws.Range("M68").Value = -2970.2222
$ws.Range("H71").Value = 3831.4666
$ws.Range("I71").Value = 3719.2222
$ws.Range("K71").Value = 18596.111
$ws.Range("M71").Value = -14852.111
$ws.Range("H82").Value = 1562.4
$ws.Range("I82").Value = 853.7
$ws.Range("J82").Value = 2979.8
$ws.Range("K82").Value = 853.7
$ws.Range("L82").Value = 2979.8
$ws.Range("M82").Value = -492.7
$ws.Range("N82").Value = -3701.8
$ws.Range("H85").Value = 1562.4
$ws.Range("I85").Value = 853.7
$ws.Range("J85").Value = 2979.8
$ws.Range("K85").Value = 853.7
$ws.Range("L85").Value = 2979.8
$ws.Range("M85").Value = 394.3
$ws.Range("N85").Value = -5475.8
$ws.Range("H132").Value = 8218.75
$ws.Range("I132").Value = 3653.1875
$ws.Range("J132").Value = 17349.875
$ws.Range("K132").Value = 10959.5625
$ws.Range("L132").Value = 52049.625
$ws.Range("M132").Value = -8429.5625
$ws.Range("N132").Value = -57109.625
$ws.Range("H138").Value = 174500
$ws.Range("J138").Value = 174500
$ws.Range("L138").Value = 174500
$ws.Range("N138").Value = -184780
$ws.Range("H140").Value = 111996.63
$ws.Range("J140").Value = 111996.63
$ws.Range("L140").Value = 111996.63
$ws.Range("N140").Value = -122356.63
$ws.Range("H141").Value = 169329.8
$ws.Range("J141").Value = 169329.8
$ws.Range("L141").Value = 169329.8
$ws.Range("N141").Value = -179689.8
